$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:P6").Value = 4

$ws.Range("A7:F7").Value = 4
$ws.Range("G7:I7").Value = 2
$ws.Range("J7:P7").Value = 4

$ws.Range("A8:B8").Value = 4
$ws.Range("C8:D8").Value = 0
$ws.Range("E8:F8").Value = 4
$ws.Range("G8:I8").Value = 2
$ws.Range("J8:P8").Value = 4

$ws.Range("A9:C9").Value = 4
$ws.Range("D9:E9").Value = 0
$ws.Range("F9").Value = 4
$ws.Range("G9:I9").Value = 2
$ws.Range("J9:P9").Value = 4

$ws.Range("A10:D10").Value = 4
$ws.Range("E10:F10").Value = 0
$ws.Range("G10:P10").Value = 4

$ws.Range("A11:E11").Value = 4
$ws.Range("F11:G11").Value = 0
$ws.Range("H11:P11").Value = 4

$ws.Range("A12:P12").Value = 4

$ws.Range("G9").Select()
